# Weekly data refresh: a new (most-recent) price observation is inserted
# at the top of the data block (row 676), pushing all existing rows down
# by one. Mirrors Excel's normal "insert a row" behaviour, which is why
# the worksheet's $A$1:$R$713 used range grows to $A$1:$R$714.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 676 (shifts 676:713 down
# to 677:714, copying formatting - e.g. the date style on column D -
# from the row above, same as Excel's native row-insert behaviour).
$ws.Rows.Item(676).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A676").Value = 6
$ws.Range("B676").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C676").Value = "Metropolitana"
$ws.Range("D676").Value = 45267
$ws.Range("E676").Value = 13
$ws.Range("F676").Value = 100112043
$ws.Range("G676").Value = "Pepino ensalada"
$ws.Range("H676").Value = "Sin especificar"
$ws.Range("I676").Value = "Primera"
$ws.Range("J676").Value = 550
$ws.Range("K676").Value = 14000
$ws.Range("L676").Value = 15000
$ws.Range("M676").Value = 14582
$ws.Range("N676").Value = "`$/caja 60 unidades"
$ws.Range("O676").Value = "Región Metropolitana"
$ws.Range("P676").Value = 243
$ws.Range("Q676").Value = 60
$ws.Range("R676").Value = "Hortaliza"
